$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Column G: per-segment area, mirroring column E's depth-weighted pattern ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Shared formula block G4:G15, relative to row 4 ((D-prevD)*B/100)
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- H2: running total of area (Atotal) ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- J2 / K2: summary cells mirroring Atotal / Qtotal ---
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Selection, matching the post-edit saved file ---
$ws.Range("J2:K2").Select()
